$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.23060401392695
$ws.Range("C2").Value = 0.1123143690809343
$ws.Range("D2").Value = 0.4930429649347303
$ws.Range("E2").Value = 0.153012318667983
$ws.Range("G2").Value = 0.002527355200859052
$ws.Range("J2").Value = 0.05488047511336092
$ws.Range("L2").Value = 0.4910442775714046
$ws.Range("N2").Value = 1.736197621498263
$ws.Range("O2").Value = 6.411123814880199
$ws.Range("B3").Value = 2.131956074550487
$ws.Range("C3").Value = 0.100124865357202
$ws.Range("D3").Value = 0.4919802920147447
$ws.Range("E3").Value = 0.1537033091196633
$ws.Range("G3").Value = 0.002531105391039279
$ws.Range("J3").Value = 0.05500228573306654
$ws.Range("L3").Value = 0.4835916280197807
$ws.Range("N3").Value = 1.755947880334318
$ws.Range("O3").Value = 6.419591120900918
$ws.Range("B4").Value = 2.072248509538554
$ws.Range("C4").Value = 0.09259618401169689
$ws.Range("D4").Value = 0.4915409414366962
$ws.Range("E4").Value = 0.1541801141137
$ws.Range("G4").Value = 0.002533532019736297
$ws.Range("J4").Value = 0.05508224912474002
$ws.Range("L4").Value = 0.4792237629983873
$ws.Range("N4").Value = 1.768718298105721
$ws.Range("O4").Value = 6.428276561215966
$ws.Range("B5").Value = 2.048135287526634
$ws.Range("C5").Value = 0.08951714893183293
$ws.Range("D5").Value = 0.4914155538296257
$ws.Range("E5").Value = 0.154387642358456
$ws.Range("G5").Value = 0.002534552168118581
$ws.Range("J5").Value = 0.05511613927705739
$ws.Range("L5").Value = 0.4774962689683946
$ws.Range("N5").Value = 1.774084204939548
$ws.Range("O5").Value = 6.432692014362232
$ws.Range("B6").Value = 2.044144512511309
$ws.Range("C6").Value = 0.08900521382751947
$ws.Range("D6").Value = 0.4913979754109761
$ws.Range("E6").Value = 0.1544229015831835
$ws.Range("G6").Value = 0.00253472345494808
$ws.Range("J6").Value = 0.05512184561572298
$ws.Range("L6").Value = 0.4772125910986773
$ws.Range("N6").Value = 1.774984985580923
$ws.Range("O6").Value = 6.433478090798019
$ws.Range("B7").Value = 2.071922425580965
$ws.Range("C7").Value = 0.09255470362440121
$ws.Range("D7").Value = 0.491539033106946
$ws.Range("E7").Value = 0.1541828593374586
$ws.Range("G7").Value = 0.002533545651054832
$ws.Range("J7").Value = 0.05508270089209955
$ws.Range("L7").Value = 0.4792002528952395
$ws.Range("N7").Value = 1.768790009248752
$ws.Range("O7").Value = 6.428332563268526
$ws.Range("B8").Value = 2.196411822350115
$ws.Range("C8").Value = 0.1081206726061623
$ws.Range("D8").Value = 0.4926323516043567
$ws.Range("E8").Value = 0.153239678347596
$ws.Range("G8").Value = 0.002528622589771169
$ws.Range("J8").Value = 0.05492140459017492
$ws.Range("L8").Value = 0.4884314828194931
$ws.Range("N8").Value = 1.742873915823591
$ws.Range("O8").Value = 6.413319321444988
$ws.Range("B9").Value = 2.447343215124647
$ws.Range("C9").Value = 0.1382913163039916
$ws.Range("D9").Value = 0.4964656572306865
$ws.Range("E9").Value = 0.1518062712295318
$ws.Range("G9").Value = 0.002519947891971097
$ws.Range("J9").Value = 0.05464595197955902
$ws.Range("L9").Value = 0.5081812093876721
$ws.Range("N9").Value = 1.6971592763724
$ws.Range("O9").Value = 6.411578450434348
$ws.Range("B10").Value = 2.635824974339698
$ws.Range("C10").Value = 0.1602400363737786
$ws.Range("D10").Value = 0.5003103975120382
$ws.Range("E10").Value = 0.1510060273706681
$ws.Range("G10").Value = 0.002514165439716449
$ws.Range("J10").Value = 0.05446823107030685
$ws.Range("L10").Value = 0.5236926422612385
$ws.Range("N10").Value = 1.666684324793156
$ws.Range("O10").Value = 6.427247523763924
$ws.Range("B11").Value = 2.722460933144816
$ws.Range("C11").Value = 0.1701775859526151
$ws.Range("D11").Value = 0.5022825542915399
$ws.Range("E11").Value = 0.1506967193242179
$ws.Range("G11").Value = 0.002511661823993345
$ws.Range("J11").Value = 0.05439268217186388
$ws.Range("L11").Value = 0.5309661551898586
$ws.Range("N11").Value = 1.653495858948375
$ws.Range("O11").Value = 6.438070057695995
$ws.Range("B12").Value = 2.755395535255616
$ws.Range("C12").Value = 0.173933846823445
$ws.Range("D12").Value = 0.5030614132539881
$ws.Range("E12").Value = 0.1505874479903824
$ws.Range("G12").Value = 0.002510731910536146
$ws.Range("J12").Value = 0.05436483139578918
$ws.Range("L12").Value = 0.5337516080072078
$ws.Range("N12").Value = 1.648598784104838
$ws.Range("O12").Value = 6.442700532253127
$ws.Range("B13").Value = 2.748296830380468
$ws.Range("C13").Value = 0.1731251762817294
$ws.Range("D13").Value = 0.5028922474941879
$ws.Range("E13").Value = 0.1506106322903609
$ws.Range("G13").Value = 0.002510931378091415
$ws.Range("J13").Value = 0.05437079590848537
$ws.Range("L13").Value = 0.5331503281315264
$ws.Range("N13").Value = 1.649649136360249
$ws.Range("O13").Value = 6.441679591207674
$ws.Range("B14").Value = 2.725167935376192
$ws.Range("C14").Value = 0.1704867541777162
$ws.Range("D14").Value = 0.5023459894921274
$ws.Range("E14").Value = 0.1506875721141405
$ws.Range("G14").Value = 0.002511584956191202
$ws.Range("J14").Value = 0.05439037570272376
$ws.Range("L14").Value = 0.5311946927991471
$ws.Range("N14").Value = 1.6530910269054
$ws.Range("O14").Value = 6.438440337816417
$ws.Range("B15").Value = 2.711017380779481
$ws.Range("C15").Value = 0.1688697460693902
$ws.Range("D15").Value = 0.5020155627072143
$ws.Range("E15").Value = 0.1507357228038941
$ws.Range("G15").Value = 0.002511987653047609
$ws.Range("J15").Value = 0.05440246748811095
$ws.Range("L15").Value = 0.5300008605488387
$ws.Range("N15").Value = 1.655211933542926
$ws.Range("O15").Value = 6.436525539813942
$ws.Range("B16").Value = 2.630181010203728
$ws.Range("C16").Value = 0.1595896382745536
$ws.Range("D16").Value = 0.5001859992286342
$ws.Range("E16").Value = 0.1510273415301349
$ws.Range("G16").Value = 0.002514331601840041
$ws.Range("J16").Value = 0.05447327463851437
$ws.Range("L16").Value = 0.523221663637969
$ws.Range("N16").Value = 1.667559806627935
$ws.Range("O16").Value = 6.426614670298704
$ws.Range("B17").Value = 2.580818797900804
$ws.Range("C17").Value = 0.1538844677604914
$ws.Range("D17").Value = 0.4991207470749401
$ws.Range("E17").Value = 0.1512202479177329
$ws.Range("G17").Value = 0.002515801963919761
$ws.Range("J17").Value = 0.05451806657886227
$ws.Range("L17").Value = 0.5191184162145106
$ws.Range("N17").Value = 1.675307702364531
$ws.Range("O17").Value = 6.421481610734133
$ws.Range("B18").Value = 2.552511240203785
$ws.Range("C18").Value = 0.1505985841227755
$ws.Range("D18").Value = 0.4985290502102657
$ws.Range("E18").Value = 0.1513363544762569
$ws.Range("G18").Value = 0.002516659622452673
$ws.Range("J18").Value = 0.05454432854354874
$ws.Range("L18").Value = 0.5167787953869549
$ws.Range("N18").Value = 1.679827591626864
$ws.Range("O18").Value = 6.418876906462003
$ws.Range("B19").Value = 2.542941299867607
$ws.Range("C19").Value = 0.1494852845336254
$ws.Range("D19").Value = 0.4983323213422182
$ws.Range("E19").Value = 0.1513765513982204
$ws.Range("G19").Value = 0.00251695206534349
$ws.Range("J19").Value = 0.05455330618283494
$ws.Range("L19").Value = 0.5159901572426975
$ws.Range("N19").Value = 1.681368852986733
$ws.Range("O19").Value = 6.418054687085316
$ws.Range("B20").Value = 2.586064775641603
$ws.Range("C20").Value = 0.1544922512019014
$ws.Range("D20").Value = 0.499231971064475
$ws.Range("E20").Value = 0.1511991795934033
$ws.Range("G20").Value = 0.002515644205766647
$ws.Range("J20").Value = 0.05451324680108538
$ws.Range("L20").Value = 0.5195530973527696
$ws.Range("N20").Value = 1.674476352557733
$ws.Range("O20").Value = 6.421992041692022
$ws.Range("B21").Value = 2.731958003244472
$ws.Range("C21").Value = 0.1712619101249118
$ws.Range("D21").Value = 0.5025055693439384
$ws.Range("E21").Value = 0.1506647598892403
$ws.Range("G21").Value = 0.002511392492679952
$ws.Range("J21").Value = 0.05438460410207391
$ws.Range("L21").Value = 0.5317682663081342
$ws.Range("N21").Value = 1.652077423960147
$ws.Range("O21").Value = 6.439377333343202
$ws.Range("B22").Value = 2.828049953853679
$ws.Range("C22").Value = 0.1821817371576628
$ws.Range("D22").Value = 0.5048318233184119
$ws.Range("E22").Value = 0.1503612760113633
$ws.Range("G22").Value = 0.002508719510315016
$ws.Range("J22").Value = 0.0543049448498083
$ws.Range("L22").Value = 0.5399329905731491
$ws.Range("N22").Value = 1.638004509135303
$ws.Range("O22").Value = 6.453842230847101
$ws.Range("B23").Value = 2.776696347663574
$ws.Range("C23").Value = 0.1763573266363778
$ws.Range("D23").Value = 0.5035731831235637
$ws.Range("E23").Value = 0.1505190654218929
$ws.Range("G23").Value = 0.00251013648394912
$ws.Range("J23").Value = 0.05434705765307069
$ws.Range("L23").Value = 0.5355587622946985
$ws.Range("N23").Value = 1.645463663122808
$ws.Range("O23").Value = 6.445837841590674
$ws.Range("B24").Value = 2.583692846490578
$ws.Range("C24").Value = 0.1542174906956575
$ws.Range("D24").Value = 0.4991816221216681
$ws.Range("E24").Value = 0.151208688377956
$ws.Range("G24").Value = 0.002515715489818559
$ws.Range("J24").Value = 0.05451542423263511
$ws.Range("L24").Value = 0.5193565176028727
$ws.Range("N24").Value = 1.67485200170734
$ws.Range("O24").Value = 6.421760196997866
$ws.Range("B25").Value = 2.378733697039308
$ws.Range("C25").Value = 0.1301674570422335
$ws.Range("D25").Value = 0.4952479455920695
$ws.Range("E25").Value = 0.1521495783107554
$ws.Range("G25").Value = 0.002522190421091079
$ws.Range("J25").Value = 0.05471612184114694
$ws.Range("L25").Value = 0.5026623371403076
$ws.Range("N25").Value = 1.708979774838273
$ws.Range("O25").Value = 6.409077617275472